$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) down from row 24 into the new row 25 first, so the
# subsequent value/formula writes keep the same number formats as the rest
# of the table (date style on B, 2-decimal % style on G/H/I/N, integer style
# on J/K).
$ws.Range("B24").Copy($ws.Range("B25"))
$ws.Range("G24").Copy($ws.Range("G25"))
$ws.Range("H24").Copy($ws.Range("H25"))
$ws.Range("I24").Copy($ws.Range("I25"))
$ws.Range("J24").Copy($ws.Range("J25"))
$ws.Range("K24").Copy($ws.Range("K25"))
$ws.Range("N24").Copy($ws.Range("N25"))

# New day row: Index 24, 2025-03-26 (Wed), start/end balances, withdrawals
# and market totals for the day.
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 45742
$ws.Range("C25").Value = "Wed"
$ws.Range("D25").Value = 9113
$ws.Range("E25").Value = 9163
$ws.Range("F25").Formula = "=E25-D25"
$ws.Range("G25").Formula = "=(E25-`$D`$2)/A25"
$ws.Range("H25").Formula = "=(E25/D25-1)*100"
$ws.Range("I25").Formula = "=(POWER((E25/`$D`$3),1/A25)-1)*100"
$ws.Range("J25").Formula = "=J24*1.013"
$ws.Range("K25").Formula = "=E25-J25"
$ws.Range("L25").Value = 0
$ws.Range("M25").Formula = "=L25+E25"
$ws.Range("N25").Formula = "=E25/`$D`$2*100"
$ws.Range("O25").Value = 87985.1
$ws.Range("P25").Value = 86736.1
$ws.Range("Q25").Formula = "=P25/`$O`$2*100"

# Keep the sheet's dimension/selection markers in step with the newly added
# row, matching where Excel would have left the cursor after entry.
[void]$ws.Range("P26").Select()
